# calories app complete for now
#
# This script applies the final batch of updates to the "calories" tracking
# sheet ("data"): it fills in the previously-blank "consumption" column (H)
# with computed totals (mostly 0, a few 1000's), marks a handful of rows as
# modified on 11/13/2022 (column C), and zeroes out the "excercise" (E)
# values on the days whose consumption is now accounted for, plus a small
# text/number fix on row 43.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("data")

# --- Column E ("excercise") value fixes ---
# Row 43 was stored as a text "-500" - make it a real number.
$ws.Range("E43").Value = -500

# These rows' -2700 excercise entries are superseded (now 0).
$noExcerciseRows = @(45, 52, 59, 66, 73, 80, 87, 94)
foreach ($r in $noExcerciseRows) {
    $ws.Range("E$r").Value = 0
}

# --- Column C ("modified") text date stamps ---
# Mark these rows as modified on 11/13/2022. The column stores dates as
# plain text (matching column B's formatting), so force a text number
# format while writing, then restore the default "Normal" style so no
# stray formatting is left behind on the cell.
$modifiedRows = @(45, 52, 56, 59, 66, 73, 80, 85, 86, 87, 94)
foreach ($r in $modifiedRows) {
    $cell = $ws.Range("C$r")
    $cell.NumberFormat = "@"
    $cell.Value = "11/13/2022"
    $cell.Style = "Normal"
}

# --- Column H ("consumption") updates ---
# Most rows simply get a computed consumption of 0.
$zeroConsumptionRows = @()
for ($r = 2; $r -le 94; $r++) {
    if ($r -eq 36 -or $r -eq 56 -or $r -eq 85 -or $r -eq 86 -or $r -eq 87) {
        continue
    }
    $zeroConsumptionRows += $r
}
foreach ($r in $zeroConsumptionRows) {
    $ws.Range("H$r").Value = 0
}

# A few rows actually consumed 1000.
$thousandConsumptionRows = @(56, 85, 87)
foreach ($r in $thousandConsumptionRows) {
    $ws.Range("H$r").Value = 1000
}

# Row 86's consumption of 1000 was entered as text, not a number - keep it
# that way (same text-preserving technique as the column C updates above).
$cell = $ws.Range("H86")
$cell.NumberFormat = "@"
$cell.Value = "1000"
$cell.Style = "Normal"
